$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(94).Insert()

$ws.Range("A94").Value = 3
$ws.Range("B94").Value = "Femacal de La Calera"
$ws.Range("C94").Value = "Coquimbo"
$ws.Range("D94").Value = 44413
$ws.Range("E94").Value = 5
$ws.Range("F94").Value = 100112037
$ws.Range("G94").Value = "Cebollín"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 310
$ws.Range("K94").Value = 3800
$ws.Range("L94").Value = 4000
$ws.Range("M94").Value = 3903
$ws.Range("N94").Value = "$/paquete 36 unidades"
$ws.Range("O94").Value = "Provincia de Quillota"
$ws.Range("P94").Value = 108
$ws.Range("Q94").Value = 36
$ws.Range("R94").Value = "Hortaliza"
